# Auto-generated Excel COM-interop script to apply the profit-table refresh
# across the ALC/ARM/BSM/CRP/GSM/LTW/WVR sheets, per the scheduled-runner diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 891.5294
$ws.Range("I41").Value = 215
$ws.Range("J41").Value = 1492.8889
$ws.Range("K41").Value = 215
$ws.Range("L41").Value = 1492.8889
$ws.Range("M41").Value = 225
$ws.Range("N41").Value = -2372.8889
$ws.Range("H64").Value = 3249.5
$ws.Range("I64").Value = 3061.875
$ws.Range("J64").Value = 4000
$ws.Range("K64").Value = 3061.875
$ws.Range("L64").Value = 4000
$ws.Range("M64").Value = -2813.875
$ws.Range("N64").Value = -4496
$ws.Range("H67").Value = 3249.5
$ws.Range("I67").Value = 3061.875
$ws.Range("J67").Value = 4000
$ws.Range("K67").Value = 3061.875
$ws.Range("L67").Value = 4000
$ws.Range("M67").Value = -2203.875
$ws.Range("N67").Value = -5716
$ws.Range("H69").Value = 6647.0586
$ws.Range("J69").Value = 3937.5
$ws.Range("L69").Value = 11812.5
$ws.Range("N69").Value = -13560.5
$ws.Range("H72").Value = 6647.0586
$ws.Range("J72").Value = 3937.5
$ws.Range("L72").Value = 35437.5
$ws.Range("N72").Value = -44173.5
$ws.Range("H76").Value = 35717212
$ws.Range("I76").Value = 38464420
$ws.Range("J76").Value = 3500
$ws.Range("K76").Value = 38464420
$ws.Range("L76").Value = 3500
$ws.Range("M76").Value = -38464105
$ws.Range("N76").Value = -4130
$ws.Range("H79").Value = 35717212
$ws.Range("I79").Value = 38464420
$ws.Range("J79").Value = 3500
$ws.Range("K79").Value = 38464420
$ws.Range("L79").Value = 3500
$ws.Range("M79").Value = -38463328
$ws.Range("N79").Value = -5684
$ws.Range("H132").Value = 6538554
$ws.Range("I132").Value = 2887
$ws.Range("J132").Value = 55556056
$ws.Range("K132").Value = 8661
$ws.Range("L132").Value = 166668168
$ws.Range("M132").Value = -6131
$ws.Range("N132").Value = -166673228
$ws.Range("H135").Value = 25000644
$ws.Range("I135").Value = 650.86664
$ws.Range("J135").Value = 100000620
$ws.Range("K135").Value = 5857.79976
$ws.Range("L135").Value = 900005580
$ws.Range("M135").Value = -3322.79976
$ws.Range("N135").Value = -900010650

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4080.58
$ws.Range("I32").Value = 3327.9512
$ws.Range("K32").Value = 3327.9512
$ws.Range("M32").Value = -3040.9512
$ws.Range("H74").Value = 55273924
$ws.Range("I74").Value = 50715196
$ws.Range("J74").Value = 66670744
$ws.Range("K74").Value = 50715196
$ws.Range("L74").Value = 66670744
$ws.Range("M74").Value = -50714322
$ws.Range("N74").Value = -66672492
$ws.Range("H77").Value = 55273924
$ws.Range("I77").Value = 50715196
$ws.Range("J77").Value = 66670744
$ws.Range("K77").Value = 253575980
$ws.Range("L77").Value = 333353720
$ws.Range("M77").Value = -253571612
$ws.Range("N77").Value = -333362456

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 451
$ws.Range("I22").Value = 451
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 451
$ws.Range("L22").Value = 0
$ws.Range("M22").ClearContents()
$ws.Range("N22").Value = -278
$ws.Range("H86").Value = 1866.19
$ws.Range("I86").Value = 1867.4694
$ws.Range("J86").Value = 1803.5
$ws.Range("K86").Value = 1867.4694
$ws.Range("L86").Value = 1803.5
$ws.Range("M86").Value = -744.4694
$ws.Range("N86").Value = -4049.5
$ws.Range("H89").Value = 1866.19
$ws.Range("I89").Value = 1867.4694
$ws.Range("J89").Value = 1803.5
$ws.Range("K89").Value = 9337.347
$ws.Range("L89").Value = 9017.5
$ws.Range("M89").Value = -3721.347
$ws.Range("N89").Value = -20249.5
$ws.Range("H94").Value = 1670.9678
$ws.Range("I94").Value = 1231.0476
$ws.Range("J94").Value = 2594.8
$ws.Range("K94").Value = 1231.0476
$ws.Range("L94").Value = 2594.8
$ws.Range("M94").Value = -780.0476000000001
$ws.Range("N94").Value = -3496.8

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2425389.8
$ws.Range("I31").Value = 1668370.5
$ws.Range("J31").Value = 3476805.5
$ws.Range("K31").Value = 1668370.5
$ws.Range("L31").Value = 3476805.5
$ws.Range("M31").Value = -1668075.5
$ws.Range("N31").Value = -3477395.5
$ws.Range("H34").Value = 2425389.8
$ws.Range("I34").Value = 1668370.5
$ws.Range("J34").Value = 3476805.5
$ws.Range("K34").Value = 1668370.5
$ws.Range("L34").Value = 3476805.5
$ws.Range("M34").Value = -1668168.5
$ws.Range("N34").Value = -3477209.5
$ws.Range("H132").Value = 1959.5952
$ws.Range("I132").Value = 1580.4615
$ws.Range("J132").Value = 2575.6875
$ws.Range("K132").Value = 4741.3845
$ws.Range("L132").Value = 7727.0625
$ws.Range("M132").Value = -2211.3845
$ws.Range("N132").Value = -12787.0625

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 17941.334
$ws.Range("I15").Value = 20000
$ws.Range("J15").Value = 17684
$ws.Range("K15").Value = 20000
$ws.Range("L15").Value = 17684
$ws.Range("M15").Value = -19712
$ws.Range("N15").Value = -18260
$ws.Range("H34").Value = 25000
$ws.Range("J34").Value = 25000
$ws.Range("L34").Value = 25000
$ws.Range("N34").Value = -25536
$ws.Range("H62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("L62").ClearContents()
$ws.Range("N62").Value = 0
$ws.Range("H65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("L65").ClearContents()
$ws.Range("N65").Value = 0
$ws.Range("H68").Value = 22000
$ws.Range("J68").Value = 22000
$ws.Range("L68").Value = 22000
$ws.Range("N68").Value = -23622
$ws.Range("H69").Value = 29000
$ws.Range("J69").Value = 29000
$ws.Range("L69").Value = 29000
$ws.Range("N69").Value = -30498
$ws.Range("H71").Value = 22000
$ws.Range("J71").Value = 22000
$ws.Range("L71").Value = 66000
$ws.Range("N71").Value = -74112
$ws.Range("H72").Value = 29000
$ws.Range("J72").Value = 29000
$ws.Range("L72").Value = 87000
$ws.Range("N72").Value = -94488
$ws.Range("H76").Value = 25000
$ws.Range("J76").Value = 25000
$ws.Range("L76").Value = 25000
$ws.Range("N76").Value = -25630
$ws.Range("H79").Value = 25000
$ws.Range("J79").Value = 25000
$ws.Range("L79").Value = 25000
$ws.Range("N79").Value = -27184
$ws.Range("H80").Value = 13737.85
$ws.Range("I80").Value = 5452.273
$ws.Range("J80").Value = 23864.666
$ws.Range("K80").Value = 5452.273
$ws.Range("L80").Value = 23864.666
$ws.Range("M80").Value = -4454.273
$ws.Range("N80").Value = -25860.666
$ws.Range("H81").Value = 17941.334
$ws.Range("I81").Value = 20000
$ws.Range("J81").Value = 17684
$ws.Range("K81").Value = 20000
$ws.Range("L81").Value = 17684
$ws.Range("M81").Value = -19002
$ws.Range("N81").Value = -19680
$ws.Range("H82").Value = 30000
$ws.Range("J82").Value = 30000
$ws.Range("L82").Value = 30000
$ws.Range("N82").Value = -30766
$ws.Range("H83").Value = 13737.85
$ws.Range("I83").Value = 5452.273
$ws.Range("J83").Value = 23864.666
$ws.Range("K83").Value = 27261.365
$ws.Range("L83").Value = 119323.33
$ws.Range("M83").Value = -22269.365
$ws.Range("N83").Value = -129307.33
$ws.Range("H84").Value = 17941.334
$ws.Range("I84").Value = 20000
$ws.Range("J84").Value = 17684
$ws.Range("K84").Value = 60000
$ws.Range("L84").Value = 53052
$ws.Range("M84").Value = -55008
$ws.Range("N84").Value = -63036
$ws.Range("H85").Value = 30000
$ws.Range("J85").Value = 30000
$ws.Range("L85").Value = 30000
$ws.Range("N85").Value = -32652
$ws.Range("H86").Value = 25511.777
$ws.Range("J86").Value = 25511.777
$ws.Range("L86").Value = 25511.777
$ws.Range("N86").Value = -27883.777
$ws.Range("H87").Value = 29333.334
$ws.Range("J87").Value = 29333.334
$ws.Range("L87").Value = 29333.334
$ws.Range("N87").Value = -31829.334
$ws.Range("H88").Value = 0
$ws.Range("J88").Value = 0
$ws.Range("L88").ClearContents()
$ws.Range("N88").Value = 0
$ws.Range("H89").Value = 25511.777
$ws.Range("J89").Value = 25511.777
$ws.Range("L89").Value = 76535.33099999999
$ws.Range("N89").Value = -88391.33099999999
$ws.Range("H90").Value = 29333.334
$ws.Range("J90").Value = 29333.334
$ws.Range("L90").Value = 88000.00199999999
$ws.Range("N90").Value = -100480.002
$ws.Range("H91").Value = 0
$ws.Range("J91").Value = 0
$ws.Range("L91").ClearContents()
$ws.Range("N91").Value = 0

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 7143155.5
$ws.Range("I55").Value = 12500243
$ws.Range("J55").Value = 371.46667
$ws.Range("K55").Value = 12500243
$ws.Range("L55").Value = 371.46667
$ws.Range("M55").Value = -12500070
$ws.Range("N55").Value = -717.46667

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 620875
$ws.Range("I132").Value = 2613.9697
$ws.Range("K132").Value = 7841.909100000001
$ws.Range("M132").Value = -5311.909100000001

Write-Output "Applied profit-table updates to ALC, ARM, BSM, CRP, GSM, LTW, WVR"
